# Updates the 'cryptos' price list (columns B:E, rows 2-51) to the new
# snapshot values from the commit, including the THORChain / FirstDigitalUSD /
# ordi / FraxShare / Cronos / Aave row re-ordering at the bottom of the table.
#
# Every Price/Volume cell is stored as literal text in the workbook (e.g.
# "43.477.26", "0.0940", "  -0.67%  "), so each write forces the cell to
# Text format before assigning the value (otherwise Excel's automatic type
# inference would silently turn values such as "6.00" or "0.100" into the
# numbers 6 / 0.1) and clears formatting back to the sheet's default style
# afterwards so no stray number-format gets left behind on the cell.

$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

$updates = @(
    @('D2', '43.477.26'),
    @('E2', '  -0.67%  '),
    @('D3', '2.291.08'),
    @('E3', '  -0.01%  '),
    @('E4', '  -0.19%  '),
    @('D5', '113.45'),
    @('E5', '  +0.51%  '),
    @('D6', '265.39'),
    @('E6', '  -1.35%  '),
    @('D7', '0.646'),
    @('E7', '  +3.44%  '),
    @('E8', '  +0.34%  '),
    @('D9', '0.614'),
    @('E9', '  -0.92%  '),
    @('D10', '47.03'),
    @('E10', '  -2.45%  '),
    @('D11', '0.0940'),
    @('E11', '  -0.92%  '),
    @('D12', '9.23'),
    @('E12', '  +1.89%  '),
    @('E13', '  +1.41%  '),
    @('D14', '15.39'),
    @('E14', '  -2.67%  '),
    @('D15', '2.623.29'),
    @('E15', '  -0.44%  '),
    @('D16', '0.871'),
    @('E16', '  +2.57%  '),
    @('D17', '2.293.54'),
    @('E17', '  +0.45%  '),
    @('D18', '43.461.67'),
    @('E18', '  -0.43%  '),
    @('D19', '0.0000110'),
    @('E19', '  +0.87%  '),
    @('D20', '6.77'),
    @('E20', '  +0.59%  '),
    @('D21', '72.46'),
    @('E21', '  +0.34%  '),
    @('D22', '2.43'),
    @('E22', '  -0.93%  '),
    @('D23', '237.25'),
    @('E23', '  +2.10%  '),
    @('D24', '2.87'),
    @('E24', '  +2.63%  '),
    @('D25', '9.42'),
    @('E25', '  -4.36%  '),
    @('E26', '  +1.75%  '),
    @('D27', '11.59'),
    @('E27', '  -0.35%  '),
    @('D28', '41.39'),
    @('E28', '  -0.79%  '),
    @('D29', '3.38'),
    @('E29', '  -0.80%  '),
    @('E30', '  -1.23%  '),
    @('D31', '173.74'),
    @('E31', '  -0.89%  '),
    @('D32', '21.85'),
    @('E32', '  +1.63%  '),
    @('D33', '0.0907'),
    @('E33', '  -1.65%  '),
    @('D34', '5.67'),
    @('E34', '  +0.35%  '),
    @('E35', '  +2.51%  '),
    @('D36', '0.0381'),
    @('E36', '  +5.00%  '),
    @('D37', '4.72'),
    @('E37', '  +0.89%  '),
    @('D38', '3.87'),
    @('E38', '  +0.45%  '),
    @('D39', '0.105'),
    @('E39', '  -2.08%  '),
    @('D40', '2.58'),
    @('E40', '  +7.55%  '),
    @('D41', '14.55'),
    @('E41', '  +6.63%  '),
    @('D42', '74.25'),
    @('E42', '  +1.45%  '),
    @('D43', '0.236'),
    @('E43', '  -2.01%  '),
    @('B44', 'THORChain'),
    @('C44', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'),
    @('D44', '6.00'),
    @('E44', '  -5.04%  '),
    @('B45', 'FirstDigitalUSD'),
    @('C45', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'),
    @('D45', '0.999'),
    @('E45', '  -0.31%  '),
    @('D46', '1.37'),
    @('E46', '  -0.64%  '),
    @('B47', 'ordi'),
    @('C47', 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'),
    @('D47', '73.89'),
    @('E47', '  +36.07%  '),
    @('D48', '1.27'),
    @('E48', '  +4.03%  '),
    @('B49', 'FraxShare'),
    @('C49', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'),
    @('D49', '8.61'),
    @('E49', '  -1.56%  '),
    @('B50', 'Cronos'),
    @('C50', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'),
    @('D50', '0.100'),
    @('E50', '  +0.85%  '),
    @('B51', 'Aave'),
    @('C51', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'),
    @('D51', '100.53'),
    @('E51', '  -2.31%  ')
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $text = $update[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.ClearFormats()
}
